$d = $word.ActiveDocument

# --- Change 1: remove the "Meta description: ..." paragraph that
# immediately follows the title heading paragraph. ---
$pMeta = $d.Paragraphs.Item(2)
if ($pMeta.Range.Text.TrimEnd() -notlike "Meta description:*") {
    throw "Unexpected paragraph 2 content: $($pMeta.Range.Text)"
}
$pMeta.Range.Delete()

# --- Change 2: before the final paragraph (the italic image-prompt
# paragraph), insert a new bold paragraph repeating the page title. ---
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($n)
$pLast.Range.InsertParagraphBefore()

$pNew = $d.Paragraphs.Item($n)
$pNewRange = $pNew.Range
$pNewRange.MoveEnd(1, -1) | Out-Null

$newParaXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Ankh of Anubis Free | Review of Play’N’Go's Online Slot Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$pNewRange.InsertXML($newParaXml)

# --- Change 3: replace the text of the (now last) italic paragraph
# with the meta-description copy, keeping its italic formatting. ---
$n2 = $d.Paragraphs.Count
$pImg = $d.Paragraphs.Item($n2)
$imgRange = $pImg.Range
$imgRange.MoveEnd(1, -1) | Out-Null
$imgRange.Text = "Discover all about Ankh of Anubis, an online slot game from Play’N’Go, with an Ancient Egypt theme centered on the god Anubis. Play it free and read our review."
